$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing component coordinates / rotation ---

# C3: Y adjusted
$ws.Range("C4").Value = -134.11199999999999

# C9: X/Y adjusted
$ws.Range("B10").Value = 105.41
$ws.Range("C10").Value = -117.348

# D1: X/Y adjusted and rotation changed 180 -> 270
$ws.Range("B11").Value = 114.64
$ws.Range("C11").Value = -128.72
$ws.Range("E11").Value = 270

# J1: X adjusted, and duplicate manual-Y note added in column K
$ws.Range("B15").Value = 144.79
$ws.Range("K15").Value = -129.69999999999999

# J2: Y adjusted
$ws.Range("C16").Value = -136.398

# R1: apply two-decimal number format to the Mid Y cell
$ws.Range("C17").NumberFormat = "0.00"

# --- Insert new component R3 (level shifter resistor) before R4 ---
$ws.Rows(19).Insert()
$ws.Range("A19").Value = "R3"
$ws.Range("B19").Value = 121.19
$ws.Range("C19").Value = -132.31
$ws.Range("D19").Value = "top"
$ws.Range("E19").Value = 90

# R5 moved to new position (now on row 21 after the insert above)
$ws.Range("B21").Value = 119.634
$ws.Range("C21").Value = -125.98399999999999

# U1 moved to new position (now on row 31 after the insert above)
$ws.Range("B31").Value = 112.05
$ws.Range("C31").Value = -136.22

# --- Append new component Q1 (level shifter transistor) at the end ---
$ws.Range("A34").Value = "Q1"
$ws.Range("B34").Value = 119.8
$ws.Range("C34").Value = -128.87
$ws.Range("D34").Value = "top"
$ws.Range("E34").Value = 270

# --- Selection moved by the author while editing ---
$ws.Range("E12").Select()
